# Reorder the comma-separated mention lists in column B so that the
# shortest / base mention comes first, matching the new shared-string text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "ink, indelible visible ink"
$ws.Range("B8").Value = "country 's elections, recent elections, many elections, these elections, elections, upcoming parliamentary elections"
$ws.Range("B10").Value = "drive, petition drive"
$ws.Range("B13").Value = "use, improper use"
